{"js": "// Office.js (Word JavaScript API) script.\n// Applies two changes described by the diff:\n//  1. Colors the \"Set up MySQL database and YAML file\" list item blue (0070C0).\n//  2. After the last paragraph (\"Resume Udemy courses...\"), appends:\n//       - an empty paragraph with a dotted bottom border\n//       - a plain empty paragraph\n//       - a new paragraph of text about the \"Adjusted approach\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- Change 1: color the \"Set up MySQL database and YAML file\" paragraph ----\nconst target = paragraphs.items.find(\n  (p) => p.text === \"Set up MySQL database and YAML file\"\n);\nif (target) {\n  target.font.color = \"#0070C0\";\n}\nawait context.sync();\n\n// ---- Change 2: append new paragraphs after the last paragraph in the body ----\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert three new paragraphs after the last one, forcing \"Normal\" style so\n// they don't inherit the preceding list-item's ListParagraph/numbering\n// formatting.\nconst borderPara = lastParagraph.insertParagraph(\"\", \"After\");\nborderPara.style = \"Normal\";\nconst blankPara = borderPara.insertParagraph(\"\", \"After\");\nblankPara.style = \"Normal\";\nconst textPara = blankPara.insertParagraph(\"\", \"After\");\ntextPara.style = \"Normal\";\nawait context.sync();\n\nfunction wrapOoxml(innerBodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    innerBodyXml +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Paragraph with just a dotted bottom paragraph border, no text.\nborderPara\n  .getRange()\n  .insertOoxml(\n    wrapOoxml(\n      '<w:p><w:pPr><w:pBdr><w:bottom w:val=\"dotted\" w:sz=\"24\" w:space=\"1\" w:color=\"auto\"/></w:pBdr></w:pPr></w:p>'\n    ),\n    \"Replace\"\n  );\nawait context.sync();\n\n// Completely empty paragraph.\nblankPara.getRange().insertOoxml(wrapOoxml(\"<w:p/>\"), \"Replace\");\nawait context.sync();\n\n// Paragraph with the new \"Adjusted approach\" text, split across two runs\n// exactly like the source edit.\ntextPara\n  .getRange()\n  .insertOoxml(\n    wrapOoxml(\n      \"<w:p>\" +\n        \"<w:r><w:t>Adjusted approach: build web pages with HTML &amp; Bootstrap with high-level detail and hard-coded features. THEN connect to back-end</w:t></w:r>\" +\n        '<w:r><w:t xml:space=\"preserve\"> one feature at a time</w:t></w:r>' +\n        \"</w:p>\"\n    ),\n    \"Replace\"\n  );\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies two changes described by the diff:\n#  1. Colors the \"Set up MySQL database and YAML file\" list item blue (0070C0).\n#  2. After the last paragraph (\"Resume Udemy courses...\"), appends:\n#       - an empty paragraph with a dotted bottom border\n#       - a plain empty paragraph\n#       - a new paragraph of text about the \"Adjusted approach\"\n\n$d = $word.ActiveDocument\n\n# ---- Change 1: color the \"Set up MySQL database and YAML file\" paragraph ----\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Set up MySQL database and YAML file\") {\n        # 0070C0 (RGB) -> Word COM color longs are BGR-ordered.\n        $red = 0x00\n        $green = 0x70\n        $blue = 0xC0\n        $p.Range.Font.Color = $blue * 65536 + $green * 256 + $red\n        break\n    }\n}\n\n# ---- Change 2: append new paragraphs after the last paragraph in the body ----\n# Insert right before the document's final (sentinel) paragraph mark so the\n# new paragraphs land after \"Resume Udemy courses...\" without disturbing it.\n$endPos = $d.Content.End - 1\n$insertRange = $d.Range($endPos, $endPos)\n\n$innerXml = '<w:p><w:pPr><w:pBdr><w:bottom w:val=\"dotted\" w:sz=\"24\" w:space=\"1\" w:color=\"auto\"/></w:pBdr></w:pPr></w:p>' + `\n            '<w:p/>' + `\n            '<w:p><w:r><w:t>Adjusted approach: build web pages with HTML &amp; Bootstrap with high-level detail and hard-coded features. THEN connect to back-end</w:t></w:r>' + `\n            '<w:r><w:t xml:space=\"preserve\"> one feature at a time</w:t></w:r></w:p>'\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:body>' + $innerXml + '</w:body></w:document>' + `\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n[void]$insertRange.InsertXML($flatOpc)\n"}
